$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 464; this shifts rows 464:529 down to 465:530
# and carries formatting (e.g. date style on column D) along.
$ws.Rows("464:464").Insert()

# Populate the newly inserted row 464 with the new record's data.
$ws.Cells.Item(464, 1).Value2 = 7
$ws.Cells.Item(464, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(464, 3).Value = "Ñuble"
$ws.Cells.Item(464, 4).Value2 = 44776
$ws.Cells.Item(464, 5).Value2 = 16
$ws.Cells.Item(464, 6).Value = "Fruta"
$ws.Cells.Item(464, 7).Value2 = 100106
$ws.Cells.Item(464, 8).Value = "Oleaginosos"
$ws.Cells.Item(464, 9).Value2 = 100106002
$ws.Cells.Item(464, 10).Value = "Palta"
$ws.Cells.Item(464, 11).Value = "Hass"
$ws.Cells.Item(464, 12).Value = "Primera"
$ws.Cells.Item(464, 13).Value2 = 120
$ws.Cells.Item(464, 14).Value2 = 26000
$ws.Cells.Item(464, 15).Value2 = 27000
$ws.Cells.Item(464, 16).Value2 = 26500
$ws.Cells.Item(464, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(464, 18).Value = "Perú"
$ws.Cells.Item(464, 19).Value2 = 2650
$ws.Cells.Item(464, 20).Value2 = 10
